$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4724.4893  # ALC!H40
$ws.Cells.Item(40, 9).Value = 4776.5  # ALC!I40
$ws.Cells.Item(40, 10).Value = 3961.6667  # ALC!J40
$ws.Cells.Item(40, 11).Value = 4776.5  # ALC!K40
$ws.Cells.Item(40, 12).Value = 3961.6667  # ALC!L40
$ws.Cells.Item(40, 13).Value = -4601.5  # ALC!M40
$ws.Cells.Item(40, 14).Value = -4311.6667  # ALC!N40

$ws.Cells.Item(107, 8).Value = 21742228  # ALC!H107
$ws.Cells.Item(107, 9).Value = 41669204  # ALC!I107
$ws.Cells.Item(107, 11).Value = 41669204  # ALC!K107
$ws.Cells.Item(107, 13).Value = -41667284  # ALC!M107

$ws.Cells.Item(113, 8).Value = 1899  # ALC!H113
$ws.Cells.Item(113, 10).Value = 0  # ALC!J113
$ws.Cells.Item(113, 12).Value = 0  # ALC!L113
$ws.Cells.Item(113, 14).ClearContents()  # ALC!N113

$ws.Cells.Item(116, 8).Value = 7946.263  # ALC!H116
$ws.Cells.Item(116, 9).Value = 6125.2856  # ALC!I116
$ws.Cells.Item(116, 10).Value = 9008.5  # ALC!J116
$ws.Cells.Item(116, 11).Value = 6125.2856  # ALC!K116
$ws.Cells.Item(116, 12).Value = 9008.5  # ALC!L116
$ws.Cells.Item(116, 13).Value = -2683.2856  # ALC!M116
$ws.Cells.Item(116, 14).Value = -15892.5  # ALC!N116

$ws.Cells.Item(135, 8).Value = 2201.9048  # ALC!H135
$ws.Cells.Item(135, 9).Value = 2052.1052  # ALC!I135
$ws.Cells.Item(135, 11).Value = 18468.9468  # ALC!K135
$ws.Cells.Item(135, 13).Value = -15933.9468  # ALC!M135

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4091.32  # ARM!H61
$ws.Cells.Item(61, 10).Value = 8299.200000000001  # ARM!J61
$ws.Cells.Item(61, 12).Value = 8299.200000000001  # ARM!L61
$ws.Cells.Item(61, 14).Value = -8723.200000000001  # ARM!N61

$ws.Cells.Item(74, 8).Value = 5530.766  # ARM!H74
$ws.Cells.Item(74, 9).Value = 5211.0293  # ARM!I74
$ws.Cells.Item(74, 11).Value = 5211.0293  # ARM!K74
$ws.Cells.Item(74, 13).Value = -4337.0293  # ARM!M74

$ws.Cells.Item(77, 8).Value = 5530.766  # ARM!H77
$ws.Cells.Item(77, 9).Value = 5211.0293  # ARM!I77
$ws.Cells.Item(77, 11).Value = 26055.1465  # ARM!K77
$ws.Cells.Item(77, 13).Value = -21687.1465  # ARM!M77

$ws.Cells.Item(122, 8).Value = 2896.923  # ARM!H122
$ws.Cells.Item(122, 9).Value = 2229.5  # ARM!I122
$ws.Cells.Item(122, 10).Value = 3469  # ARM!J122
$ws.Cells.Item(122, 11).Value = 6688.5  # ARM!K122
$ws.Cells.Item(122, 12).Value = 10407  # ARM!L122
$ws.Cells.Item(122, 13).Value = -4238.5  # ARM!M122
$ws.Cells.Item(122, 14).Value = -15307  # ARM!N122

$ws.Cells.Item(132, 8).Value = 2887.0454  # ARM!H132
$ws.Cells.Item(132, 9).Value = 2857.861  # ARM!I132
$ws.Cells.Item(132, 10).Value = 3018.375  # ARM!J132
$ws.Cells.Item(132, 11).Value = 8573.582999999999  # ARM!K132
$ws.Cells.Item(132, 12).Value = 9055.125  # ARM!L132
$ws.Cells.Item(132, 13).Value = -6043.582999999999  # ARM!M132
$ws.Cells.Item(132, 14).Value = -14115.125  # ARM!N132

$ws.Cells.Item(136, 8).Value = 4091.32  # ARM!H136
$ws.Cells.Item(136, 10).Value = 8299.200000000001  # ARM!J136
$ws.Cells.Item(136, 12).Value = 24897.6  # ARM!L136
$ws.Cells.Item(136, 14).Value = -29997.6  # ARM!N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 158.55556  # BSM!H11
$ws.Cells.Item(11, 10).Value = 141.14285  # BSM!J11
$ws.Cells.Item(11, 12).Value = 141.14285  # BSM!L11
$ws.Cells.Item(11, 14).Value = -421.14285  # BSM!N11

$ws.Cells.Item(107, 8).Value = 4391.8438  # BSM!H107
$ws.Cells.Item(107, 9).Value = 4392.485  # BSM!I107
$ws.Cells.Item(107, 11).Value = 4392.485  # BSM!K107
$ws.Cells.Item(107, 13).Value = -2472.485  # BSM!M107

$ws.Cells.Item(134, 8).Value = 2091.1843  # BSM!H134
$ws.Cells.Item(134, 9).Value = 1931.4865  # BSM!I134
$ws.Cells.Item(134, 11).Value = 5794.4595  # BSM!K134
$ws.Cells.Item(134, 13).Value = -3259.4595  # BSM!M134

$ws.Cells.Item(139, 8).Value = 144997  # BSM!H139
$ws.Cells.Item(139, 10).Value = 144997  # BSM!J139
$ws.Cells.Item(139, 12).Value = 144997  # BSM!L139
$ws.Cells.Item(139, 14).Value = -155277  # BSM!N139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3539.0356  # CRP!H31
$ws.Cells.Item(31, 9).Value = 3113.6428  # CRP!I31
$ws.Cells.Item(31, 10).Value = 3964.4285  # CRP!J31
$ws.Cells.Item(31, 11).Value = 3113.6428  # CRP!K31
$ws.Cells.Item(31, 12).Value = 3964.4285  # CRP!L31
$ws.Cells.Item(31, 13).Value = -2818.6428  # CRP!M31
$ws.Cells.Item(31, 14).Value = -4554.4285  # CRP!N31

$ws.Cells.Item(34, 8).Value = 3539.0356  # CRP!H34
$ws.Cells.Item(34, 9).Value = 3113.6428  # CRP!I34
$ws.Cells.Item(34, 10).Value = 3964.4285  # CRP!J34
$ws.Cells.Item(34, 11).Value = 3113.6428  # CRP!K34
$ws.Cells.Item(34, 12).Value = 3964.4285  # CRP!L34
$ws.Cells.Item(34, 13).Value = -2911.6428  # CRP!M34
$ws.Cells.Item(34, 14).Value = -4368.4285  # CRP!N34

$ws.Cells.Item(58, 8).Value = 8175.1787  # CRP!H58
$ws.Cells.Item(58, 9).Value = 9477.471  # CRP!I58
$ws.Cells.Item(58, 11).Value = 9477.471  # CRP!K58
$ws.Cells.Item(58, 13).Value = -9274.471  # CRP!M58

$ws.Cells.Item(132, 8).Value = 3852.4  # CRP!H132
$ws.Cells.Item(132, 9).Value = 3664.4595  # CRP!I132
$ws.Cells.Item(132, 11).Value = 10993.3785  # CRP!K132
$ws.Cells.Item(132, 13).Value = -8463.378499999999  # CRP!M132

$ws.Cells.Item(136, 8).Value = 8175.1787  # CRP!H136
$ws.Cells.Item(136, 9).Value = 9477.471  # CRP!I136
$ws.Cells.Item(136, 11).Value = 28432.413  # CRP!K136
$ws.Cells.Item(136, 13).Value = -25882.413  # CRP!M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 97.35714  # CUL!H2
$ws.Cells.Item(2, 9).Value = 34.666668  # CUL!I2
$ws.Cells.Item(2, 10).Value = 144.375  # CUL!J2
$ws.Cells.Item(2, 11).Value = 208.000008  # CUL!K2
$ws.Cells.Item(2, 12).Value = 866.25  # CUL!L2
$ws.Cells.Item(2, 13).Value = -95.00000800000001  # CUL!M2
$ws.Cells.Item(2, 14).Value = -1092.25  # CUL!N2

$ws.Cells.Item(12, 8).Value = 116.92308  # CUL!H12
$ws.Cells.Item(12, 10).Value = 129.2  # CUL!J12
$ws.Cells.Item(12, 12).Value = 387.6  # CUL!L12
$ws.Cells.Item(12, 14).Value = -733.5999999999999  # CUL!N12

$ws.Cells.Item(68, 8).Value = 2186.5454  # CUL!H68
$ws.Cells.Item(68, 10).Value = 2215.2  # CUL!J68
$ws.Cells.Item(68, 12).Value = 6645.599999999999  # CUL!L68
$ws.Cells.Item(68, 14).Value = -8267.599999999999  # CUL!N68

$ws.Cells.Item(71, 8).Value = 2186.5454  # CUL!H71
$ws.Cells.Item(71, 10).Value = 2215.2  # CUL!J71
$ws.Cells.Item(71, 12).Value = 19936.8  # CUL!L71
$ws.Cells.Item(71, 14).Value = -28048.8  # CUL!N71

$ws.Cells.Item(107, 8).Value = 373.16666  # CUL!H107
$ws.Cells.Item(107, 10).Value = 403.4  # CUL!J107
$ws.Cells.Item(107, 12).Value = 1210.2  # CUL!L107
$ws.Cells.Item(107, 14).Value = -5050.2  # CUL!N107

$ws.Cells.Item(116, 8).Value = 8500943  # CUL!H116
$ws.Cells.Item(116, 9).Value = 8500943  # CUL!I116
$ws.Cells.Item(116, 10).Value = 0  # CUL!J116
$ws.Cells.Item(116, 11).Value = 25502829  # CUL!K116
$ws.Cells.Item(116, 12).Value = 0  # CUL!L116
$ws.Cells.Item(116, 13).Value = -25499387  # CUL!M116
$ws.Cells.Item(116, 14).ClearContents()  # CUL!N116

$ws.Cells.Item(132, 8).Value = 2786.6667  # CUL!H132
$ws.Cells.Item(132, 9).Value = 2523.1538  # CUL!I132
$ws.Cells.Item(132, 11).Value = 22708.3842  # CUL!K132
$ws.Cells.Item(132, 13).Value = -20178.3842  # CUL!M132

$ws.Cells.Item(137, 8).Value = 11460.125  # CUL!H137
$ws.Cells.Item(137, 9).Value = 7798.8  # CUL!I137
$ws.Cells.Item(137, 10).Value = 12423.632  # CUL!J137
$ws.Cells.Item(137, 11).Value = 23396.4  # CUL!K137
$ws.Cells.Item(137, 12).Value = 37270.896  # CUL!L137
$ws.Cells.Item(137, 13).Value = -18296.4  # CUL!M137
$ws.Cells.Item(137, 14).Value = -47470.896  # CUL!N137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 32.444443  # GSM!H2
$ws.Cells.Item(2, 9).Value = 27.571428  # GSM!I2
$ws.Cells.Item(2, 11).Value = 27.571428  # GSM!K2
$ws.Cells.Item(2, 13).Value = 85.428572  # GSM!M2

$ws.Cells.Item(41, 8).Value = 22285.285  # GSM!H41
$ws.Cells.Item(41, 9).Value = 24666.166  # GSM!I41
$ws.Cells.Item(41, 11).Value = 24666.166  # GSM!K41
$ws.Cells.Item(41, 13).Value = -24311.166  # GSM!M41

$ws.Cells.Item(80, 8).Value = 12003648  # GSM!H80
$ws.Cells.Item(80, 9).Value = 30002136  # GSM!I80
$ws.Cells.Item(80, 11).Value = 30002136  # GSM!K80
$ws.Cells.Item(80, 13).Value = -30001138  # GSM!M80

$ws.Cells.Item(83, 8).Value = 12003648  # GSM!H83
$ws.Cells.Item(83, 9).Value = 30002136  # GSM!I83
$ws.Cells.Item(83, 11).Value = 150010680  # GSM!K83
$ws.Cells.Item(83, 13).Value = -150005688  # GSM!M83

$ws.Cells.Item(126, 8).Value = 4939.9443  # GSM!H126
$ws.Cells.Item(126, 9).Value = 4642.5  # GSM!I126
$ws.Cells.Item(126, 10).Value = 5007.5454  # GSM!J126
$ws.Cells.Item(126, 11).Value = 13927.5  # GSM!K126
$ws.Cells.Item(126, 12).Value = 15022.6362  # GSM!L126
$ws.Cells.Item(126, 13).Value = -11457.5  # GSM!M126
$ws.Cells.Item(126, 14).Value = -19962.6362  # GSM!N126

$ws.Cells.Item(132, 8).Value = 5150.3125  # GSM!H132
$ws.Cells.Item(132, 9).Value = 5335.524  # GSM!I132
$ws.Cells.Item(132, 10).Value = 4796.727  # GSM!J132
$ws.Cells.Item(132, 11).Value = 16006.572  # GSM!K132
$ws.Cells.Item(132, 12).Value = 14390.181  # GSM!L132
$ws.Cells.Item(132, 13).Value = -13476.572  # GSM!M132
$ws.Cells.Item(132, 14).Value = -19450.181  # GSM!N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5001.773  # LTW!H7
$ws.Cells.Item(7, 9).Value = 5066.1665  # LTW!I7
$ws.Cells.Item(7, 10).Value = 4977.625  # LTW!J7
$ws.Cells.Item(7, 11).Value = 5066.1665  # LTW!K7
$ws.Cells.Item(7, 12).Value = 4977.625  # LTW!L7
$ws.Cells.Item(7, 13).Value = -4954.1665  # LTW!M7
$ws.Cells.Item(7, 14).Value = -5201.625  # LTW!N7

$ws.Cells.Item(46, 8).Value = 3029.25  # LTW!H46
$ws.Cells.Item(46, 10).Value = 3000.1333  # LTW!J46
$ws.Cells.Item(46, 12).Value = 3000.1333  # LTW!L46
$ws.Cells.Item(46, 14).Value = -3376.1333  # LTW!N46

$ws.Cells.Item(61, 8).Value = 9678  # LTW!H61
$ws.Cells.Item(61, 10).Value = 2301  # LTW!J61
$ws.Cells.Item(61, 12).Value = 2301  # LTW!L61
$ws.Cells.Item(61, 14).Value = -2705  # LTW!N61

$ws.Cells.Item(82, 8).Value = 379.85  # LTW!H82
$ws.Cells.Item(82, 9).Value = 385.28125  # LTW!I82
$ws.Cells.Item(82, 10).Value = 249.5  # LTW!J82
$ws.Cells.Item(82, 11).Value = 385.28125  # LTW!K82
$ws.Cells.Item(82, 12).Value = 249.5  # LTW!L82
$ws.Cells.Item(82, 13).Value = -24.28125  # LTW!M82
$ws.Cells.Item(82, 14).Value = -971.5  # LTW!N82

$ws.Cells.Item(85, 8).Value = 379.85  # LTW!H85
$ws.Cells.Item(85, 9).Value = 385.28125  # LTW!I85
$ws.Cells.Item(85, 10).Value = 249.5  # LTW!J85
$ws.Cells.Item(85, 11).Value = 385.28125  # LTW!K85
$ws.Cells.Item(85, 12).Value = 249.5  # LTW!L85
$ws.Cells.Item(85, 13).Value = 862.71875  # LTW!M85
$ws.Cells.Item(85, 14).Value = -2745.5  # LTW!N85

$ws.Cells.Item(113, 8).Value = 9678  # LTW!H113
$ws.Cells.Item(113, 10).Value = 2301  # LTW!J113
$ws.Cells.Item(113, 12).Value = 2301  # LTW!L113
$ws.Cells.Item(113, 14).Value = -6641  # LTW!N113

$ws.Cells.Item(126, 8).Value = 5001.773  # LTW!H126
$ws.Cells.Item(126, 9).Value = 5066.1665  # LTW!I126
$ws.Cells.Item(126, 10).Value = 4977.625  # LTW!J126
$ws.Cells.Item(126, 11).Value = 15198.4995  # LTW!K126
$ws.Cells.Item(126, 12).Value = 14932.875  # LTW!L126
$ws.Cells.Item(126, 13).Value = -12728.4995  # LTW!M126
$ws.Cells.Item(126, 14).Value = -19872.875  # LTW!N126

$ws.Cells.Item(136, 8).Value = 4747128.5  # LTW!H136
$ws.Cells.Item(136, 9).Value = 11259865  # LTW!I136
$ws.Cells.Item(136, 10).Value = 10592.682  # LTW!J136
$ws.Cells.Item(136, 11).Value = 33779595  # LTW!K136
$ws.Cells.Item(136, 12).Value = 31778.046  # LTW!L136
$ws.Cells.Item(136, 13).Value = -33777045  # LTW!M136
$ws.Cells.Item(136, 14).Value = -36878.046  # LTW!N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5351510.5  # WVR!H81
$ws.Cells.Item(81, 9).Value = 5684730  # WVR!I81
$ws.Cells.Item(81, 11).Value = 11369460  # WVR!K81
$ws.Cells.Item(81, 13).Value = -11368399  # WVR!M81

$ws.Cells.Item(84, 8).Value = 5351510.5  # WVR!H84
$ws.Cells.Item(84, 9).Value = 5684730  # WVR!I84
$ws.Cells.Item(84, 11).Value = 56847300  # WVR!K84
$ws.Cells.Item(84, 13).Value = -56841996  # WVR!M84

$ws.Cells.Item(132, 8).Value = 3996.1936  # WVR!H132
$ws.Cells.Item(132, 9).Value = 4178.0356  # WVR!I132
$ws.Cells.Item(132, 11).Value = 12534.1068  # WVR!K132
$ws.Cells.Item(132, 13).Value = -10004.1068  # WVR!M132
